$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 49 - this shifts the existing data rows
# (old rows 49..132) down by one (new rows 50..133), growing the used
# range from A1:R132 to A1:R133, matching the diff.
$ws.Rows.Item(49).Insert()

# Populate the newly-inserted row 49 with the latest week's record
# (same Mercado/Categoria/Variedad/Unidad/Origen/Clasificacion block as
# the surrounding rows; only the date + market figures differ).
$ws.Range("A49").Value = 10
$ws.Range("B49").Value = "Vega Modelo de Temuco"
$ws.Range("C49").Value = "La Araucanía"
$ws.Range("D49").Value = 44469
$ws.Range("E49").Value = 9
$ws.Range("F49").Value = 100112005
$ws.Range("G49").Value = "Puerro"
$ws.Range("H49").Value = "Azul de Maquehue"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 140
$ws.Range("K49").Value = 6000
$ws.Range("L49").Value = 7000
$ws.Range("M49").Value = 6571
$ws.Range("N49").Value = "$/docena de paquetes"
$ws.Range("O49").Value = "Provincia de Cautín"
$ws.Range("P49").Value = 548
$ws.Range("Q49").Value = 12
$ws.Range("R49").Value = "Hortaliza"
